$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for Price cells whose new values look numeric,
# so they remain text (matching the original text-formatted price column)
# instead of being auto-coerced into floating point numbers by Excel.
$textCells = @("D5", "D7", "D8", "D11", "D15", "D16", "D18", "D19", "D22", "D23", "D25", "D26", "D28", "D30", "D31", "D33", "D38", "D39", "D40", "D44", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value2 = "27.924.44"
$ws.Range("E2").Value2 = "  +0.10%  "
$ws.Range("D3").Value2 = "1.636.36"
$ws.Range("E3").Value2 = "  +0.26%  "
$ws.Range("E4").Value2 = "  +0.00%  "
$ws.Range("D5").Value2 = "211.84"
$ws.Range("E5").Value2 = "  +0.05%  "
$ws.Range("D7").Value2 = "0.999"
$ws.Range("E7").Value2 = "  -0.02%  "
$ws.Range("D8").Value2 = "23.38"
$ws.Range("E8").Value2 = "  +0.71%  "
$ws.Range("E9").Value2 = "  -0.34%  "
$ws.Range("E10").Value2 = "  -0.17%  "
$ws.Range("D11").Value2 = "0.0884"
$ws.Range("E11").Value2 = "  +0.54%  "
$ws.Range("D12").Value2 = "1.867.44"
$ws.Range("E12").Value2 = "  +0.20%  "
$ws.Range("D13").Value2 = "1.638.88"
$ws.Range("E13").Value2 = "  +0.38%  "
$ws.Range("E14").Value2 = "  -0.92%  "
$ws.Range("D15").Value2 = "0.563"
$ws.Range("E15").Value2 = "  -0.28%  "
$ws.Range("D16").Value2 = "65.31"
$ws.Range("E16").Value2 = "  +0.15%  "
$ws.Range("D17").Value2 = "27.924.15"
$ws.Range("E17").Value2 = "  +0.09%  "
$ws.Range("D18").Value2 = "229.15"
$ws.Range("E18").Value2 = "  -0.34%  "
$ws.Range("D19").Value2 = "7.73"
$ws.Range("E19").Value2 = "  +3.00%  "
$ws.Range("E20").Value2 = "  -0.16%  "
$ws.Range("D22").Value2 = "4.37"
$ws.Range("E22").Value2 = "  +0.24%  "
$ws.Range("D23").Value2 = "10.16"
$ws.Range("E23").Value2 = "  -1.85%  "
$ws.Range("E24").Value2 = "  +0.67%  "
$ws.Range("D25").Value2 = "155.98"
$ws.Range("E25").Value2 = "  +1.44%  "
$ws.Range("D26").Value2 = "6.96"
$ws.Range("E26").Value2 = "  +0.12%  "
$ws.Range("E27").Value2 = "  +0.05%  "
$ws.Range("D28").Value2 = "15.56"
$ws.Range("E28").Value2 = "  -0.31%  "
$ws.Range("E29").Value2 = "  -0.05%  "
$ws.Range("D30").Value2 = "1.18"
$ws.Range("E30").Value2 = "  +0.15%  "
$ws.Range("D31").Value2 = "0.0482"
$ws.Range("E31").Value2 = "  -0.01%  "
$ws.Range("E32").Value2 = "  +1.16%  "
$ws.Range("D33").Value2 = "3.11"
$ws.Range("E33").Value2 = "  +1.36%  "
$ws.Range("D34").Value2 = "1.400.15"
$ws.Range("E34").Value2 = "  +0.17%  "
$ws.Range("E35").Value2 = "  +3.19%  "
$ws.Range("E36").Value2 = "  +0.75%  "
$ws.Range("E37").Value2 = "  -0.69%  "
$ws.Range("D38").Value2 = "0.0171"
$ws.Range("E38").Value2 = "  +0.28%  "
$ws.Range("D39").Value2 = "0.559"
$ws.Range("E39").Value2 = "  -0.09%  "
$ws.Range("D40").Value2 = "0.852"
$ws.Range("E40").Value2 = "  -2.06%  "
$ws.Range("E41").Value2 = "  +0.05%  "
$ws.Range("E42").Value2 = "  -1.13%  "
$ws.Range("E43").Value2 = "  +2.80%  "
$ws.Range("D44").Value2 = "66.10"
$ws.Range("E44").Value2 = "  -1.13%  "
$ws.Range("E45").Value2 = "  -1.18%  "
$ws.Range("D46").Value2 = "1.775.73"
$ws.Range("E46").Value2 = "  +0.08%  "
$ws.Range("D48").Value2 = "88.74"
$ws.Range("E48").Value2 = "  +1.29%  "
$ws.Range("D49").Value2 = "0.103"
$ws.Range("E49").Value2 = "  +2.18%  "
$ws.Range("E50").Value2 = "  -0.41%  "
$ws.Range("D51").Value2 = "7.64"
$ws.Range("E51").Value2 = "  +2.00%  "
